# Update the 'Last Updated' timestamp on the Metadata sheet
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 07:11 PM"

# Insert a new 'Top Losers' sheet right after 'Top Gainers' (mirrors the
# existing Top Gainers layout: Icon / Stock / Latest / Weekly / Monthly)
$gainers = $wb.Worksheets.Item("Top Gainers")
$losers = $wb.Worksheets.Add($null, $gainers)
$losers.Name = "Top Losers"

# Match the outline/page setup conventions used by the rest of the workbook
$losers.Outline.SummaryRow = 1
$losers.Outline.SummaryColumn = 1
$losers.PageSetup.LeftMargin = 54
$losers.PageSetup.RightMargin = 54
$losers.PageSetup.TopMargin = 72
$losers.PageSetup.BottomMargin = 72
$losers.PageSetup.HeaderMargin = 36
$losers.PageSetup.FooterMargin = 36

# Header row, styled the same way as the other data sheets' header row
$header = $losers.Range("A1:E1")
$losers.Cells.Item(1,1).Value = "Icon"
$losers.Cells.Item(1,2).Value = "Stock"
$losers.Cells.Item(1,3).Value = "Latest"
$losers.Cells.Item(1,4).Value = "Weekly"
$losers.Cells.Item(1,5).Value = "Monthly"
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data rows: full list of Top Losers (latest % change, weekly, monthly)
$arr = New-Object 'object[,]' 75,5
$arr[0,0] = "📉"
$arr[0,1] = "IIFLCAPS"
$arr[0,2] = -9.721500000000001
$arr[0,3] = -5.9431
$arr[0,4] = 20.4574
$arr[1,0] = "📉"
$arr[1,1] = "COHANCE"
$arr[1,2] = -8.5204
$arr[1,3] = -10.2636
$arr[1,4] = -10.5392
$arr[2,0] = "📉"
$arr[2,1] = "KHAICHEM"
$arr[2,2] = -8.4635
$arr[2,3] = -6.7205
$arr[2,4] = -4.6455
$arr[3,0] = "📉"
$arr[3,1] = "MOTILALOFS"
$arr[3,2] = -7.6514
$arr[3,3] = -2.7408
$arr[3,4] = 12.7798
$arr[4,0] = "📉"
$arr[4,1] = "VERANDA"
$arr[4,2] = -7.1035
$arr[4,3] = -7.4096
$arr[4,4] = 8.066700000000001
$arr[5,0] = "📉"
$arr[5,1] = "MOLDTKPAC"
$arr[5,2] = -7.0505
$arr[5,3] = -1.2462
$arr[5,4] = 0.466
$arr[6,0] = "📉"
$arr[6,1] = "NETWEB"
$arr[6,2] = -6.7177
$arr[6,3] = 5.5715
$arr[6,4] = 7.5565
$arr[7,0] = "📉"
$arr[7,1] = "FABTECH"
$arr[7,2] = -6.147
$arr[7,3] = 15.6232
$arr[7,4] = "N/A"
$arr[8,0] = "📉"
$arr[8,1] = "CARTRADE"
$arr[8,2] = -5.7253
$arr[8,3] = 16.1193
$arr[8,4] = 20.6359
$arr[9,0] = "📉"
$arr[9,1] = "TVSELECT"
$arr[9,2] = -5.6153
$arr[9,3] = -0.9738
$arr[9,4] = -2.9968
$arr[10,0] = "📉"
$arr[10,1] = "NSLNISP"
$arr[10,2] = -5.4542
$arr[10,3] = 1.3037
$arr[10,4] = 0.4681
$arr[11,0] = "📉"
$arr[11,1] = "NAM-INDIA"
$arr[11,2] = -5.1992
$arr[11,3] = -7.0279
$arr[11,4] = -1.2034
$arr[12,0] = "📉"
$arr[12,1] = "KICL"
$arr[12,2] = -5.0437
$arr[12,3] = -0.8405
$arr[12,4] = 21.7969
$arr[13,0] = "📉"
$arr[13,1] = "CCCL"
$arr[13,2] = -5.0146
$arr[13,3] = -4.576
$arr[13,4] = -12.4759
$arr[14,0] = "📉"
$arr[14,1] = "CREDITACC"
$arr[14,2] = -4.9692
$arr[14,3] = -1.3216
$arr[14,4] = 3.7319
$arr[15,0] = "📉"
$arr[15,1] = "KALAMANDIR"
$arr[15,2] = -4.8415
$arr[15,3] = 1.7451
$arr[15,4] = 25.9996
$arr[16,0] = "📉"
$arr[16,1] = "CRAMC"
$arr[16,2] = -4.7668
$arr[16,3] = 5.978
$arr[16,4] = "N/A"
$arr[17,0] = "📉"
$arr[17,1] = "SMLISUZU"
$arr[17,2] = -4.7654
$arr[17,3] = 4.993
$arr[17,4] = -2.8236
$arr[18,0] = "📉"
$arr[18,1] = "MANAKCOAT"
$arr[18,2] = -4.5883
$arr[18,3] = -6.011
$arr[18,4] = 24.8093
$arr[19,0] = "📉"
$arr[19,1] = "HDFCAMC"
$arr[19,2] = -4.401
$arr[19,3] = -2.6247
$arr[19,4] = -2.4311
$arr[20,0] = "📉"
$arr[20,1] = "ATHERENERG"
$arr[20,2] = -4.0945
$arr[20,3] = -0.0142
$arr[20,4] = 24.8806
$arr[21,0] = "📉"
$arr[21,1] = "SHAREINDIA"
$arr[21,2] = -4.0806
$arr[21,3] = -1.6889
$arr[21,4] = 54.7217
$arr[22,0] = "📉"
$arr[22,1] = "INDIQUBE"
$arr[22,2] = -4.0805
$arr[22,3] = -4.8408
$arr[22,4] = -6.7982
$arr[23,0] = "📉"
$arr[23,1] = "CHENNPETRO"
$arr[23,2] = -3.973
$arr[23,3] = 4.8359
$arr[23,4] = 6.8488
$arr[24,0] = "📉"
$arr[24,1] = "KFINTECH"
$arr[24,2] = -3.873
$arr[24,3] = -1.7487
$arr[24,4] = 7.3785
$arr[25,0] = "📉"
$arr[25,1] = "360ONE"
$arr[25,2] = -3.8488
$arr[25,3] = -4.976
$arr[25,4] = 10.0293
$arr[26,0] = "📉"
$arr[26,1] = "SMSPHARMA"
$arr[26,2] = -3.7339
$arr[26,3] = -3.0871
$arr[26,4] = 17.4387
$arr[27,0] = "📉"
$arr[27,1] = "BHARATWIRE"
$arr[27,2] = -3.5327
$arr[27,3] = 22.8336
$arr[27,4] = 23.8979
$arr[28,0] = "📉"
$arr[28,1] = "ABSLAMC"
$arr[28,2] = -3.5313
$arr[28,3] = -5.9355
$arr[28,4] = -1.2887
$arr[29,0] = "📉"
$arr[29,1] = "SUMMITSEC"
$arr[29,2] = -3.4113
$arr[29,3] = -1.6476
$arr[29,4] = 6.0097
$arr[30,0] = "📉"
$arr[30,1] = "SPLPETRO"
$arr[30,2] = -3.3984
$arr[30,3] = -5.0241
$arr[30,4] = -7.769
$arr[31,0] = "📉"
$arr[31,1] = "MPSLTD"
$arr[31,2] = -3.2782
$arr[31,3] = -4.6314
$arr[31,4] = 2.1755
$arr[32,0] = "📉"
$arr[32,1] = "CAMS"
$arr[32,2] = -3.2545
$arr[32,3] = -0.6366000000000001
$arr[32,4] = 2.5781
$arr[33,0] = "📉"
$arr[33,1] = "SPARC"
$arr[33,2] = -3.1709
$arr[33,3] = 4.8337
$arr[33,4] = 6.3311
$arr[34,0] = "📉"
$arr[34,1] = "PRUDENT"
$arr[34,2] = -3.127
$arr[34,3] = -3.5103
$arr[34,4] = 2.1213
$arr[35,0] = "📉"
$arr[35,1] = "ANANDRATHI"
$arr[35,2] = -3.0775
$arr[35,3] = -0.8672
$arr[35,4] = 9.1835
$arr[36,0] = "📉"
$arr[36,1] = "NLCINDIA"
$arr[36,2] = -3.0757
$arr[36,3] = -4.5618
$arr[36,4] = -11.6431
$arr[37,0] = "📉"
$arr[37,1] = "YATRA"
$arr[37,2] = -3.0403
$arr[37,3] = -2.8455
$arr[37,4] = 7.3711
$arr[38,0] = "📉"
$arr[38,1] = "BOSCHLTD"
$arr[38,2] = -3.0099
$arr[38,3] = -3.123
$arr[38,4] = -2.0055
$arr[39,0] = "📉"
$arr[39,1] = "DRREDDY"
$arr[39,2] = -2.9859
$arr[39,3] = -2.5475
$arr[39,4] = 2.2228
$arr[40,0] = "📉"
$arr[40,1] = "ROSSTECH"
$arr[40,2] = -2.9778
$arr[40,3] = 1.9028
$arr[40,4] = -6.8057
$arr[41,0] = "📉"
$arr[41,1] = "OAL"
$arr[41,2] = -2.9496
$arr[41,3] = -1.278
$arr[41,4] = 8.7362
$arr[42,0] = "📉"
$arr[42,1] = "ENDURANCE"
$arr[42,2] = -2.939
$arr[42,3] = -2.2945
$arr[42,4] = 3.4531
$arr[43,0] = "📉"
$arr[43,1] = "POLICYBZR"
$arr[43,2] = -2.907
$arr[43,3] = 2.2365
$arr[43,4] = 1.2573
$arr[44,0] = "📉"
$arr[44,1] = "AYMSYNTEX"
$arr[44,2] = -2.9052
$arr[44,3] = -0.3705
$arr[44,4] = -10.494
$arr[45,0] = "📉"
$arr[45,1] = "DIGITIDE"
$arr[45,2] = -2.8795
$arr[45,3] = 3.2317
$arr[45,4] = 6.2968
$arr[46,0] = "📉"
$arr[46,1] = "RUBICON"
$arr[46,2] = -2.8712
$arr[46,3] = -2.9654
$arr[46,4] = "N/A"
$arr[47,0] = "📉"
$arr[47,1] = "STARHEALTH"
$arr[47,2] = -2.8707
$arr[47,3] = -1.5572
$arr[47,4] = 7.5434
$arr[48,0] = "📉"
$arr[48,1] = "KIRIINDUS"
$arr[48,2] = -2.8411
$arr[48,3] = -1.3849
$arr[48,4] = 1.4335
$arr[49,0] = "📉"
$arr[49,1] = "UNIMECH"
$arr[49,2] = -2.8008
$arr[49,3] = -1.6104
$arr[49,4] = -0.4585
$arr[50,0] = "📉"
$arr[50,1] = "TTKPRESTIG"
$arr[50,2] = -2.7438
$arr[50,3] = 8.001200000000001
$arr[50,4] = 9.650499999999999
$arr[51,0] = "📉"
$arr[51,1] = "PFOCUS"
$arr[51,2] = -2.7039
$arr[51,3] = -2.6276
$arr[51,4] = -1.2163
$arr[52,0] = "📉"
$arr[52,1] = "ALLDIGI"
$arr[52,2] = -2.6342
$arr[52,3] = -0.2306
$arr[52,4] = -5.3103
$arr[53,0] = "📉"
$arr[53,1] = "PRIVISCL"
$arr[53,2] = -2.6288
$arr[53,3] = -2.1048
$arr[53,4] = 19.7451
$arr[54,0] = "📉"
$arr[54,1] = "CANHLIFE"
$arr[54,2] = -2.6148
$arr[54,3] = 5.2953
$arr[54,4] = "N/A"
$arr[55,0] = "📉"
$arr[55,1] = "GKENERGY"
$arr[55,2] = -2.6122
$arr[55,3] = -9.807700000000001
$arr[55,4] = 23.2758
$arr[56,0] = "📉"
$arr[56,1] = "SGFIN"
$arr[56,2] = -2.592
$arr[56,3] = -0.06270000000000001
$arr[56,4] = 11.7235
$arr[57,0] = "📉"
$arr[57,1] = "ARVINDFASN"
$arr[57,2] = -2.549
$arr[57,3] = -2.9892
$arr[57,4] = -4.4223
$arr[58,0] = "📉"
$arr[58,1] = "EDELWEISS"
$arr[58,2] = -2.5422
$arr[58,3] = -3.3745
$arr[58,4] = 8.5305
$arr[59,0] = "📉"
$arr[59,1] = "SAMHI"
$arr[59,2] = -2.5284
$arr[59,3] = 1.8231
$arr[59,4] = 2.8516
$arr[60,0] = "📉"
$arr[60,1] = "TBOTEK"
$arr[60,2] = -2.524
$arr[60,3] = -3.5732
$arr[60,4] = 1.036
$arr[61,0] = "📉"
$arr[61,1] = "UJJIVANSFB"
$arr[61,2] = -2.5201
$arr[61,3] = 0.3845
$arr[61,4] = 12.6645
$arr[62,0] = "📉"
$arr[62,1] = "AMBER"
$arr[62,2] = -2.5098
$arr[62,3] = -0.1082
$arr[62,4] = 2.763
$arr[63,0] = "📉"
$arr[63,1] = "GRPLTD"
$arr[63,2] = -2.4898
$arr[63,3] = -5.9894
$arr[63,4] = -5.4586
$arr[64,0] = "📉"
$arr[64,1] = "NESCO"
$arr[64,2] = -2.4722
$arr[64,3] = 1.9934
$arr[64,4] = 3.8931
$arr[65,0] = "📉"
$arr[65,1] = "PILANIINVS"
$arr[65,2] = -2.4546
$arr[65,3] = -0.7907
$arr[65,4] = 4.267
$arr[66,0] = "📉"
$arr[66,1] = "NSIL"
$arr[66,2] = -2.4088
$arr[66,3] = -1.7646
$arr[66,4] = 4.7431
$arr[67,0] = "📉"
$arr[67,1] = "COALINDIA"
$arr[67,2] = -2.4016
$arr[67,3] = -3.058
$arr[67,4] = -2.0387
$arr[68,0] = "📉"
$arr[68,1] = "JNKINDIA"
$arr[68,2] = -2.3482
$arr[68,3] = -2.8371
$arr[68,4] = 4.2622
$arr[69,0] = "📉"
$arr[69,1] = "FCL"
$arr[69,2] = -2.3453
$arr[69,3] = -2.616
$arr[69,4] = -0.02
$arr[70,0] = "📉"
$arr[70,1] = "DEEDEV"
$arr[70,2] = -2.3334
$arr[70,3] = -6.6528
$arr[70,4] = -7.4227
$arr[71,0] = "📉"
$arr[71,1] = "WEALTH"
$arr[71,2] = -2.2793
$arr[71,3] = -3.8356
$arr[71,4] = -2.7981
$arr[72,0] = "📉"
$arr[72,1] = "RATNAMANI"
$arr[72,2] = -2.2788
$arr[72,3] = -0.4626
$arr[72,4] = 0.8712
$arr[73,0] = "📉"
$arr[73,1] = "CSBBANK"
$arr[73,2] = -2.2695
$arr[73,3] = 2.3137
$arr[73,4] = 10.6999
$arr[74,0] = "📉"
$arr[74,1] = "BBOX"
$arr[74,2] = -2.2639
$arr[74,3] = -4.7636
$arr[74,4] = 5.1528

$losers.Range("A2:E76").Value = $arr
